# Add a new weekly ranking sheet "2025-12-22" by duplicating the previous
# week (2025-12-15) and updating the title/author/latest_episode columns.
$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-12-15")

# Duplicate the previous week sheet right after itself (i.e. at the end)
# so it inherits the header style, column layout and formatting.
$src.Copy([System.Reflection.Missing]::Value, $src)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2025-12-22"

# New ranking data for the week of 2025-12-22.
# Columns: title, author, latest_episode (rank in column A is unchanged: 1..50)
$data = @(
    @('生徒会にも穴はある！', 'むちまろ', '第142話	私もぎゅっと!!'),
    @('転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～', 'zunta(作画) はらわたさいぞう(原作)', '第34話：プロのテク③'),
    @('時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―', '光永康則', '第７２話『先端停止』③'),
    @('勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが', '絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)', '第6話 後編'),
    @('蜘蛛ですが、なにか？', 'かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)', '第77話その2'),
    @('いとこのこ', 'いぬちく(著者)', '第43話'),
    @('帰ってください！ 阿久津さん', '長岡太一(著者)', '第200話'),
    @('元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～', '沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)', '第81話その2'),
    @('怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった', '菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)', '第17話後半'),
    @('実は俺、最強でした？', '原作：澄守 彩 漫画：高橋 愛', '第135話　腕試ししちゃいましょ♪・前編'),
    @('ぽんドロイド！ はまさん', 'はれやまはれぞう(著者)', '第14話'),
    @('女友達は頼めば意外とヤらせてくれる', 'ろくろ(漫画) 鏡遊(原作)', '第27話②'),
    @('異世界魔王と召喚少女の奴隷魔術', '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大', '第131話　幕間（後編）'),
    @('異世界のんびり農家', '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)', '第313話'),
    @('貞操逆転世界で頼めばヤれると噂の俺', '澄田佑貴(漫画) aaa168（スリーエー）(原作)', '第3話'),
    @('勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～', '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり', '第５５話　封印を解く器用貧乏（４）'),
    @('小林さんちのメイドラゴン', 'クール教信者', '第156話'),
    @('バキ外伝 烈海王は異世界転生しても一向にかまわんッッ', '板垣恵介 猪原賽 陸井栄史', '第86話　暇潰し'),
    @('リビルドワールド', '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)', '第76話➁'),
    @('アザミヤコを好きになる', 'ユニティコング(原作) ツノニガウ(作画)', '第12話前編'),
    @('異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～', 'モリタ Ｕ４ nima', '第14話（５）　春とおぼっちゃまとピクニックランチ（５）'),
    @('ダンジョンの幼なじみ', '久真やすひさ(著者)', '第2回人気投票結果発表！'),
    @('世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜', '戸賀 環 坂木持丸 riritto', '第58話①　おすそわけをしてみた'),
    @('貞操逆転世界の童貞辺境領主騎士', '柳瀬こたつ（漫画） 道造（原作） めろん２２（キャラクター原案）', '第11話　やむなき犠牲（後編）'),
    @('願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜', 'ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)', '第9話-2：黄昏の茶会'),
    @('よくわからないけれど異世界に転生していたようです', '内々けやき あし カオミン', '第145話 よくわからないけれどナメられているみたいです（１）'),
    @('路地裏で拾った女の子がバッドエンド後の乙女ゲームのヒロインだった件', 'カボチャマスク(原作) 樋乃えなが(作画) へいろー(キャラクター原案)', '第2話'),
    @('骸骨騎士様、只今異世界へお出掛け中', 'サワノアキラ（漫画） 秤猿鬼（原作） KeG（キャラクター原案）', '第65話　エルフ族の決断Ⅳ'),
    @('聖者無双', '漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime', '第95話　奴隷の扱い・戦闘準備（後半）'),
    @('アイドル辞めるけど結婚してくれますか!?', '三吉汐美(著者)', '第19話後半'),
    @('配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信', '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)', '第5話前編'),
    @('お気楽領主の楽しい領地防衛 ～生産系魔術で名もなき村を最強の城塞都市に～', '青色まろ（漫画） 赤池宗（原作） 転（原作イラスト）', '第36話　出陣'),
    @('くらいあの子としたいこと', '碇マナツ(著者)', '第88話'),
    @('落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～', '村上よしゆき 茨木野 あるてら', '第４４話　勇者、S級ランクに認定される。弟は、ミカエルとクエストに向かったら、ラファエルが降ってくる（１）'),
    @('『おっぱい揉みたい』って叫んだら、妹の友達と付き合うことになりました。', '凪木エコ(原作) 白クマシェイク(キャラクター原案) 逢沢もにょ(作画)', '第16話'),
    @('バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～', '板垣恵介 林たかあき', '第60話 間合いの緊張感'),
    @('理想のヒモ生活', '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)', '第90話　その3'),
    @('追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。', '六志麻あさ 業務用餅 kisui', '第７７話'),
    @('BL世界に転生したので、モブ女子とラブコメします。', 'karl(著者)', '読切'),
    @('十年目、帰還を諦めた転移者はいまさら主人公になる', '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう', '第21話②'),
    @('姫様“拷問”の時間です', '原作:春原ロビンソン　漫画:ひらけい', '拷問162'),
    @('回復術士のやり直し', '月夜涙(原作) 羽賀ソウケン(漫画) しおこんぶ(キャラクター原案)', '第75話-1'),
    @('ふかふかダンジョン攻略記～俺の異世界転生冒険譚～', 'KAKERU', '第71話「ファントム・アレイ」（後半）'),
    @('地味子な三葉さんが僕を誘惑する', 'はぶらえる(著者)', '第13話前半'),
    @('Lv２からチートだった元勇者候補のまったり異世界ライフ', '糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）', '第64話　混血児'),
    @('ゲーム悪役貴族に転生した俺は、チート筋肉で無双する', '昼行燈（原作） しいたけ元帥（漫画）', '第36話'),
    @('婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版', '漫画/すたひろ 原作/Y.A', 'chapter76【40話①】'),
    @('最弱貴族に転生したので悪役たちを集めてみた', '空野進 sorani ファルまろ', '第15話　最弱貴族、領地を成長させる（２）'),
    @('世界最高の暗殺者、異世界貴族に転生する', '月夜涙(原作) 皇ハマオ(漫画) れい亜(キャラクター原案)', '第40話-2'),
    @('黄金の経験値', '原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)', '第21話（前編）')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# Restore the originally active sheet/selection so we only change what the
# ranking update requires.
$wb.Worksheets.Item("Sheet1").Activate()

Write-Output "Added sheet 2025-12-22 with $($wb.Worksheets.Count) total worksheets"
